# Update both "wilcox_table" and "wilcox_table_selected" sheets: the rows
# describing each environmental variable are being re-ordered (table went
# from being sorted by cruise to being sorted by month). Columns A
# (Variables), G (statistic) and H (p) travel together with each variable;
# columns B-F (.y., group1, group2, n1, n2) are identical across all rows
# and stay put.

$wb = $excel.ActiveWorkbook

# New row order (name -> statistic, p) for the full table (sheet "wilcox_table")
$fullOrder = @(
    @{ Name = "CN";           Statistic = 20;   P = 0.429 },
    @{ Name = "Chla";         Statistic = 30;   P = 0.00433 },
    @{ Name = "Clay";         Statistic = 24;   P = 0.12 },
    @{ Name = "D50";          Statistic = 9;    P = 0.329 },
    @{ Name = "DRM";          Statistic = 13;   P = 0.792 },
    @{ Name = "Density";      Statistic = 20;   P = 0.429 },
    @{ Name = "Depth";        Statistic = 19.5; P = 0.462 },
    @{ Name = "Fluorescence"; Statistic = 18;   P = 0.662 },
    @{ Name = "Oxygen";       Statistic = 30;   P = 0.00433 },
    @{ Name = "Porosity";     Statistic = 6;    P = 0.126 },
    @{ Name = "Salinity";     Statistic = 20;   P = 0.429 },
    @{ Name = "Sand";         Statistic = 16;   P = 0.927 },
    @{ Name = "SigmaTheta";   Statistic = 20;   P = 0.429 },
    @{ Name = "Silt";         Statistic = 6;    P = 0.12 },
    @{ Name = "TN";           Statistic = 15;   P = 1 },
    @{ Name = "TOC";          Statistic = 18;   P = 0.662 },
    @{ Name = "Temperature";  Statistic = 10;   P = 0.429 },
    @{ Name = "Transmission"; Statistic = 27;   P = 0.0303 },
    @{ Name = "WC";           Statistic = 6;    P = 0.126 },
    @{ Name = "delta13C";     Statistic = 2;    P = 0.0173 }
)

# New row order (name -> statistic, p) for the "selected" subset table
$selectedOrder = @(
    @{ Name = "CN";           Statistic = 20; P = 0.429 },
    @{ Name = "Chla";         Statistic = 30; P = 0.00433 },
    @{ Name = "D50";          Statistic = 9;  P = 0.329 },
    @{ Name = "Fluorescence"; Statistic = 18; P = 0.662 },
    @{ Name = "Porosity";     Statistic = 6;  P = 0.126 },
    @{ Name = "TOC";          Statistic = 18; P = 0.662 },
    @{ Name = "Temperature";  Statistic = 10; P = 0.429 }
)

function Update-WilcoxSheet {
    param($Worksheet, $Order)

    $row = 2
    foreach ($item in $Order) {
        $Worksheet.Cells.Item($row, 1).Value = $item.Name
        $Worksheet.Cells.Item($row, 7).Value = $item.Statistic
        $Worksheet.Cells.Item($row, 8).Value = $item.P
        $row++
    }
}

$wsFull = $wb.Worksheets.Item("wilcox_table")
$wsSelected = $wb.Worksheets.Item("wilcox_table_selected")

Update-WilcoxSheet $wsFull $fullOrder
Update-WilcoxSheet $wsSelected $selectedOrder
